# "Generate Report for Handoff"
# Adds two new handed-off files (two .png dependency rows) to the
# localization-status report: one row per file on every sheet
# (Overview, zh-cn, de-de), refreshes the existing "source .md" row's
# handoff timestamp / hashed target filenames, and wires up the new
# hyperlinks (source file + target file) exactly like the existing rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Shared literal values re-used across sheets/rows
# ---------------------------------------------------------------------
$oldGuid       = "b5c90828-ac5c-4d04-b734-d38bacea5dcc"
$newGuid       = "40c88eae-0f71-48d7-8021-39a6c1956b1d"

$mdName        = "$newGuid.md"
$png1Name      = "9555e5bb-0bcb-45ba-a629-2e1d866ebe7e.png"
$png2Name      = "e8964933-c1dd-4979-b14c-780feeefc41c.png"

$zhHash        = "10864e6bb0d594264194a0a5c37e389c23fd7123"
$zhXlf         = "$newGuid.$zhHash.zh-cn.xlf"
$deXlf         = "$newGuid.$zhHash.de-de.xlf"

$png1Target    = "000fa002b582c91135a406dc6ad2daf1fdc34f17.png"
$png2Target    = "a24fed01c7779a8d65a8c1ad721fa68c267ee2fd.png"

$status        = "Ready for handoff"
$handoffDate   = "2016-03-22 21:07:38"
$handoffDt     = "2016-03-22 21:07:34"
$epoch         = "0001-01-01 00:00:00"
$include       = "Include"
$isDependency  = "IsDependency"
$dependencyFrom= "e2e\$mdName"

$mdUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/c247a411d3db7eb4f3ce903a3258c081a2f82239/e2e/$mdName"
$png1Url = "https://github.com/OpenLocalizationTest/oltest/blob/c247a411d3db7eb4f3ce903a3258c081a2f82239/e2e/$png1Name"
$png2Url = "https://github.com/OpenLocalizationTest/oltest/blob/c247a411d3db7eb4f3ce903a3258c081a2f82239/e2e/$png2Name"

$zhBase     = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2dbbbd789f0f4e278534373b1010978ca4c6a5f2/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht"
$deBase     = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/32f086b8f7c4aab526b9d93dba9ddcb368f6a206/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht"

$zhXlfUrl     = "$zhBase/$zhXlf"
$deXlfUrl     = "$deBase/$deXlf"
$png1TgtUrlZh = "$zhBase/$png1Target"
$png2TgtUrlZh = "$zhBase/$png2Target"
$png1TgtUrlDe = "$deBase/$png1Target"
$png2TgtUrlDe = "$deBase/$png2Target"

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# =======================================================================
# Sheet "Overview"
# =======================================================================
$ov = $wb.Worksheets.Item("Overview")

# refresh row 2 (existing .md file) in place
$ov.Range("A2").Hyperlinks.Delete()
$ov.Range("A2").Value = $mdName
$ov.Range("B2").Value = $status
$ov.Range("C2").Value = $status
$ov.Range("D2").Value = $handoffDate
$ov.Hyperlinks.Add($ov.Range("A2"), $mdUrl, "", "", $mdName) | Out-Null
$ov.Range("A2").Style = "HyperLink"
$ov.Range("D2").NumberFormat = $dateFmt

# new row 3 (first .png dependency)
$ov.Range("A3").Value = $png1Name
$ov.Range("B3").Value = $status
$ov.Range("C3").Value = $status
$ov.Range("D3").Value = $handoffDate
$ov.Hyperlinks.Add($ov.Range("A3"), $png1Url, "", "", $png1Name) | Out-Null
$ov.Range("A3").Style = "HyperLink"
$ov.Range("D3").NumberFormat = $dateFmt

# new row 4 (second .png dependency)
$ov.Range("A4").Value = $png2Name
$ov.Range("B4").Value = $status
$ov.Range("C4").Value = $status
$ov.Range("D4").Value = $handoffDate
$ov.Hyperlinks.Add($ov.Range("A4"), $png2Url, "", "", $png2Name) | Out-Null
$ov.Range("A4").Style = "HyperLink"
$ov.Range("D4").NumberFormat = $dateFmt

# =======================================================================
# Sheet "zh-cn"
# =======================================================================
$zh = $wb.Worksheets.Item("zh-cn")

# refresh row 2 (existing .md file) in place
$zh.Range("A2").Hyperlinks.Delete()
$zh.Range("D2").Hyperlinks.Delete()
$zh.Range("A2").Value = $mdName
$zh.Range("B2").Value = ".md"
$zh.Range("C2").Value = $status
$zh.Range("D2").Value = $zhXlf
$zh.Range("E2").Value = $handoffDt
$zh.Range("H2").Value = $epoch
$zh.Range("J2").Value = $include
$zh.Hyperlinks.Add($zh.Range("A2"), $mdUrl, "", "", $mdName) | Out-Null
$zh.Hyperlinks.Add($zh.Range("D2"), $zhXlfUrl, "", "", $zhXlf) | Out-Null
$zh.Range("A2").Style = "HyperLink"
$zh.Range("D2").Style = "HyperLink"

# new row 3 (first .png dependency)
$zh.Range("A3").Value = $png1Name
$zh.Range("B3").Value = ".png"
$zh.Range("C3").Value = $status
$zh.Range("D3").Value = $png1Target
$zh.Range("E3").Value = $handoffDt
$zh.Range("H3").Value = $epoch
$zh.Range("J3").Value = $isDependency
$zh.Range("K3").Value = $dependencyFrom
$zh.Hyperlinks.Add($zh.Range("A3"), $png1Url, "", "", $png1Name) | Out-Null
$zh.Hyperlinks.Add($zh.Range("D3"), $png1TgtUrlZh, "", "", $png1Target) | Out-Null
$zh.Range("A3").Style = "HyperLink"
$zh.Range("D3").Style = "HyperLink"

# new row 4 (second .png dependency)
$zh.Range("A4").Value = $png2Name
$zh.Range("B4").Value = ".png"
$zh.Range("C4").Value = $status
$zh.Range("D4").Value = $png2Target
$zh.Range("E4").Value = $handoffDt
$zh.Range("H4").Value = $epoch
$zh.Range("J4").Value = $isDependency
$zh.Range("K4").Value = $dependencyFrom
$zh.Hyperlinks.Add($zh.Range("A4"), $png2Url, "", "", $png2Name) | Out-Null
$zh.Hyperlinks.Add($zh.Range("D4"), $png2TgtUrlZh, "", "", $png2Target) | Out-Null
$zh.Range("A4").Style = "HyperLink"
$zh.Range("D4").Style = "HyperLink"

# =======================================================================
# Sheet "de-de"
# =======================================================================
$de = $wb.Worksheets.Item("de-de")

# refresh row 2 (existing .md file) in place
$de.Range("A2").Hyperlinks.Delete()
$de.Range("D2").Hyperlinks.Delete()
$de.Range("A2").Value = $mdName
$de.Range("B2").Value = ".md"
$de.Range("C2").Value = $status
$de.Range("D2").Value = $deXlf
$de.Range("E2").Value = $handoffDate
$de.Range("H2").Value = $epoch
$de.Range("J2").Value = $include
$de.Hyperlinks.Add($de.Range("A2"), $mdUrl, "", "", $mdName) | Out-Null
$de.Hyperlinks.Add($de.Range("D2"), $deXlfUrl, "", "", $deXlf) | Out-Null
$de.Range("A2").Style = "HyperLink"
$de.Range("D2").Style = "HyperLink"

# new row 3 (first .png dependency)
$de.Range("A3").Value = $png1Name
$de.Range("B3").Value = ".png"
$de.Range("C3").Value = $status
$de.Range("D3").Value = $png1Target
$de.Range("E3").Value = $handoffDate
$de.Range("H3").Value = $epoch
$de.Range("J3").Value = $isDependency
$de.Range("K3").Value = $dependencyFrom
$de.Hyperlinks.Add($de.Range("A3"), $png1Url, "", "", $png1Name) | Out-Null
$de.Hyperlinks.Add($de.Range("D3"), $png1TgtUrlDe, "", "", $png1Target) | Out-Null
$de.Range("A3").Style = "HyperLink"
$de.Range("D3").Style = "HyperLink"

# new row 4 (second .png dependency)
$de.Range("A4").Value = $png2Name
$de.Range("B4").Value = ".png"
$de.Range("C4").Value = $status
$de.Range("D4").Value = $png2Target
$de.Range("E4").Value = $handoffDate
$de.Range("H4").Value = $epoch
$de.Range("J4").Value = $isDependency
$de.Range("K4").Value = $dependencyFrom
$de.Hyperlinks.Add($de.Range("A4"), $png2Url, "", "", $png2Name) | Out-Null
$de.Hyperlinks.Add($de.Range("D4"), $png2TgtUrlDe, "", "", $png2Target) | Out-Null
$de.Range("A4").Style = "HyperLink"
$de.Range("D4").Style = "HyperLink"
